# Apply "Penalty Reward System" data shift to forecast_summary workbook.
# Sheet "Forecast Comparison": each week's Week_Start_Date (col B) moves forward
# one week, and the MyForecast value (col D) is updated.
# Sheet "Summary": several aggregate metrics are recalculated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Forecast Comparison
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Make sure column B keeps storing the week-start dates as plain text
# (as they were originally authored), not auto-converted Excel date serials.
$ws1.Range("B2:B17").NumberFormat = "@"

$weekStarts = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecasts = @(16, 15, 15, 15, 16, 16, 17, 16, 15, 15, 14, 14, 14, 14, 14, 14)

for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $weekStarts[$i]
    $ws1.Cells.Item($row, 4).Value = $myForecasts[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: Summary
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

# All "Value" entries on this sheet are stored as plain text (even the
# numeric-looking ones), so force text formatting before writing so Excel
# doesn't silently convert them into numbers or dates.
$ws2.Range("B2:B15").NumberFormat = "@"

$ws2.Range("B2").Value = "2023-01-01 to 2025-01-05"
$ws2.Range("B7").Value = "10"
$ws2.Range("B8").Value = "1524 units"
$ws2.Range("B9").Value = "240"
$ws2.Range("B10").Value = "126"
$ws2.Range("B11").Value = "61"
$ws2.Range("B12").Value = "17"
$ws2.Range("B14").Value = "14"
$ws2.Range("B15").Value = "2025-03-23"

Write-Host "Penalty Reward System edits applied"
